$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the two new columns (Mail / Phone) ---
# Copy the existing header formatting (bold, border, centered) from I1
# onto J1:K1, then fill in their text.
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Mail"
$ws.Range("K1").Value = "Phone"

# --- Column J ("Mail") ---
$ws.Range("J2").Value = " "
$ws.Range("J3").Value = " CAREERS@RTA.VN "
$ws.Range("J4").Value = " HR@CONTEMI.COM.VN`n "
$ws.Range("J5").Value = " CAREERS@VN.ZALORA.COM "
$ws.Range("J6").Value = " HTA@LRGLOBAL.COM "
$ws.Range("J7").Value = " CAREERS@VN.ZALORA.COM "
$ws.Range("J8").Value = " "
$ws.Range("J9").Value = " "
$ws.Range("J10").Value = " BCG-HCMC@BCG.COM "
$ws.Range("J11").Value = " "
$ws.Range("J12").Value = " "
$ws.Range("J13").Value = " HR@VATGIA.COM "
$ws.Range("J14").Value = " BCG-HCMC@BCG.COM "
$ws.Range("J15").Value = " "
$ws.Range("J16").Value = " PEOPLE@EPSILON-MOBILE.COM "
$ws.Range("J17").Value = " HR@EPSILON-MOBILE.COM "
$ws.Range("J18").Value = " HR@EPSILON-MOBILE.COM "
$ws.Range("J19").Value = " HR@EPSILON-MOBILE.COM "
$ws.Range("J20").Value = " HR@EPSILON-MOBILE.COM "
$ws.Range("J21").Value = " "

# --- Column K ("Phone") ---
$ws.Range("K2").Value = " "
$ws.Range("K3").Value = " "
# Phone number: force text interpretation (leading apostrophe) so the
# leading zero survives and the cell keeps quote-prefix formatting.
$ws.Range("K4").Value = "'0835471164 "
$ws.Range("K5").Value = " "
$ws.Range("K6").Value = " "
$ws.Range("K7").Value = " "
$ws.Range("K8").Value = " "
$ws.Range("K9").Value = " "
$ws.Range("K10").Value = " "
$ws.Range("K11").Value = " "
$ws.Range("K12").Value = " "
$ws.Range("K13").Value = " "
$ws.Range("K14").Value = " "
$ws.Range("K15").Value = " "
$ws.Range("K16").Value = " "
$ws.Range("K17").Value = " "
$ws.Range("K18").Value = " "
$ws.Range("K19").Value = " "
$ws.Range("K20").Value = " "
$ws.Range("K21").Value = " "

# Row 2's mail address got corrected after the rest of the sheet was
# filled in (first with a stray trailing parenthesis, then fixed).
$ws.Range("J2").Value = " CAREERS@RTA.VN APPLY@INTERNSHIP.EDU.VN)"
$ws.Range("J2").Value = " CAREERS@RTA.VN APPLY@INTERNSHIP.EDU.VN"

# --- Column widths ---
$ws.Columns.Item(9).ColumnWidth = 50.166666666666664
$ws.Columns.Item(10).ColumnWidth = 45.666666666666664
$ws.Columns.Item(11).ColumnWidth = 16.666666666666668

# --- View: zoom + selection ---
$ws.Range("J2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
